$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Fill column B on Sheet1 with squares 1..9
$squares = @(1, 4, 9, 16, 25, 36, 49, 64, 81)
for ($i = 0; $i -lt $squares.Length; $i++) {
    $row = $i + 1
    $ws1.Cells.Item($row, 2).Value = $squares[$i]
}

# Add a second sheet named "Hoja 2" after the first sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja 2"
$ws2.Range("A1").Value = "Hola"

# Keep the original sheet as the active one
$ws1.Activate()
